# Updated NATMI TPM-derived values for the Col5a1-Sdc3 ligand-receptor sheet.
# Source data (TPM per cluster) changed upstream; ligand/receptor expression,
# derived-specificity and edge-weight columns (E:J, M:T) are refreshed per row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2026263333333334
$ws.Range("H2").Value = 0.6078790000000001
$ws.Range("I2").Value = 0.001145895592151193
$ws.Range("J2").Value = 0.001145895592151193
$ws.Range("M2").Value = 24.75542533333333
$ws.Range("N2").Value = 74.26627599999999
$ws.Range("O2").Value = 0.7762421087066456
$ws.Range("P2").Value = 0.7762421087066456
$ws.Range("Q2").Value = 5.016101065400444
$ws.Range("R2").Value = 45.144909588604
$ws.Range("S2").Value = 0.0008894924108090921
$ws.Range("T2").Value = 0.000889492410809092

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2026263333333334
$ws.Range("H3").Value = 0.6078790000000001
$ws.Range("I3").Value = 0.001145895592151193
$ws.Range("J3").Value = 0.001145895592151193
$ws.Range("M3").Value = 3.818542
$ws.Range("O3").Value = 0.1197358984688377
$ws.Range("P3").Value = 0.1197358984688377
$ws.Range("Q3").Value = 0.7737371641393335
$ws.Range("R3").Value = 6.963634477254001
$ws.Range("S3").Value = 0.0001372048382777039
$ws.Range("T3").Value = 0.0001372048382777038

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2026263333333334
$ws.Range("H4").Value = 0.6078790000000001
$ws.Range("I4").Value = 0.001145895592151193
$ws.Range("J4").Value = 0.001145895592151193
$ws.Range("M4").Value = 3.317404
$ws.Range("N4").Value = 9.952211999999999
$ws.Range("O4").Value = 0.1040219928245168
$ws.Range("P4").Value = 0.1040219928245168
$ws.Range("Q4").Value = 0.6721934087053334
$ws.Range("R4").Value = 6.049740678348
$ws.Range("S4").Value = 0.0001191983430643968
$ws.Range("T4").Value = 0.0001191983430643968

# Row 5
$ws.Range("I5").Value = 0.9410232311014506
$ws.Range("J5").Value = 0.9410232311014505
$ws.Range("M5").Value = 24.75542533333333
$ws.Range("N5").Value = 74.26627599999999
$ws.Range("O5").Value = 0.7762421087066456
$ws.Range("P5").Value = 0.7762421087066456
$ws.Range("Q5").Value = 4119.282476017893
$ws.Range("R5").Value = 37073.54228416103
$ws.Range("S5").Value = 0.7304618572521311
$ws.Range("T5").Value = 0.730461857252131

# Row 6
$ws.Range("I6").Value = 0.9410232311014506
$ws.Range("J6").Value = 0.9410232311014505
$ws.Range("M6").Value = 3.818542
$ws.Range("O6").Value = 0.1197358984688377
$ws.Range("P6").Value = 0.1197358984688377
$ws.Range("Q6").Value = 635.4022575955601
$ws.Range("R6").Value = 5718.62031836004
$ws.Range("S6").Value = 0.1126742620559809
$ws.Range("T6").Value = 0.1126742620559809

# Row 7
$ws.Range("I7").Value = 0.9410232311014506
$ws.Range("J7").Value = 0.9410232311014505
$ws.Range("M7").Value = 3.317404
$ws.Range("N7").Value = 9.952211999999999
$ws.Range("O7").Value = 0.1040219928245168
$ws.Range("P7").Value = 0.1040219928245168
$ws.Range("Q7").Value = 552.01330532872
$ws.Range("R7").Value = 4968.11974795848
$ws.Range("S7").Value = 0.09788711179333867
$ws.Range("T7").Value = 0.09788711179333869

# Row 8
$ws.Range("G8").Value = 10.226113
$ws.Range("H8").Value = 30.678339
$ws.Range("I8").Value = 0.05783087330639819
$ws.Range("J8").Value = 0.05783087330639819
$ws.Range("M8").Value = 24.75542533333333
$ws.Range("N8").Value = 74.26627599999999
$ws.Range("O8").Value = 0.7762421087066456
$ws.Range("P8").Value = 0.7762421087066456
$ws.Range("Q8").Value = 253.1517768217293
$ws.Range("R8").Value = 2278.365991395564
$ws.Range("S8").Value = 0.04489075904370539
$ws.Range("T8").Value = 0.04489075904370539

# Row 9
$ws.Range("G9").Value = 10.226113
$ws.Range("H9").Value = 30.678339
$ws.Range("I9").Value = 0.05783087330639819
$ws.Range("J9").Value = 0.05783087330639819
$ws.Range("M9").Value = 3.818542
$ws.Range("O9").Value = 0.1197358984688377
$ws.Range("P9").Value = 0.1197358984688377
$ws.Range("Q9").Value = 39.048841987246
$ws.Range("R9").Value = 351.439577885214
$ws.Range("S9").Value = 0.006924431574579109
$ws.Range("T9").Value = 0.006924431574579109

# Row 10
$ws.Range("G10").Value = 10.226113
$ws.Range("H10").Value = 30.678339
$ws.Range("I10").Value = 0.05783087330639819
$ws.Range("J10").Value = 0.05783087330639819
$ws.Range("M10").Value = 3.317404
$ws.Range("N10").Value = 9.952211999999999
$ws.Range("O10").Value = 0.1040219928245168
$ws.Range("P10").Value = 0.1040219928245168
$ws.Range("Q10").Value = 33.924148170652
$ws.Range("R10").Value = 305.317333535868
$ws.Range("S10").Value = 0.006015682688113691
$ws.Range("T10").Value = 0.006015682688113692
